$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new data row (date + kilométrage reading)
$ws.Range("A26").Copy()
$ws.Range("A27").PasteSpecial(-4122)  # xlPasteFormats, keep same date style as above
$ws.Range("A27").Value = 43761
$ws.Range("B27").Value = 726

# Update the current selection to mirror the next empty cell, as Excel
# would do after data entry on row 27
$ws.Range("B28").Select()
